# Monthly labor update - January 2023
# Applies the January 2023 LFS update to the "Compiled" and "LFS Reports" sheets.

$wb = $excel.ActiveWorkbook

$compiled = $wb.Worksheets.Item("Compiled")
$lfsReports = $wb.Worksheets.Item("LFS Reports")

# ---------------------------------------------------------------------------
# 1. "Compiled" sheet - new period rows (labels first, so new shared strings are
#    created in the same order they appear in the canonical workbook).
# ---------------------------------------------------------------------------

$newPeriods = @(
    @{ Row = 260; Label = "2022 Annual"; Year = "2022"; Period = "Annual" },
    @{ Row = 261; Label = "2023 Jan";    Year = "2023"; Period = "Jan" },
    @{ Row = 262; Label = "2023 Feb";    Year = "2023"; Period = "Feb" },
    @{ Row = 263; Label = "2023 Mar";    Year = "2023"; Period = "Mar" },
    @{ Row = 264; Label = "2023 Apr";    Year = "2023"; Period = "Apr" },
    @{ Row = 265; Label = "2023 May";    Year = "2023"; Period = "May" },
    @{ Row = 266; Label = "2023 Jun";    Year = "2023"; Period = "Jun" },
    @{ Row = 267; Label = "2023 Jul";    Year = "2023"; Period = "Jul" },
    @{ Row = 268; Label = "2023 Aug";    Year = "2023"; Period = "Aug" },
    @{ Row = 269; Label = "2023 Sep";    Year = "2023"; Period = "Sep" },
    @{ Row = 270; Label = "2023 Oct";    Year = "2023"; Period = "Oct" },
    @{ Row = 271; Label = "2023 Nov";    Year = "2023"; Period = "Nov" },
    @{ Row = 272; Label = "2023 Dec";    Year = "2023"; Period = "Dec" },
    @{ Row = 273; Label = "2023 Annual"; Year = "2023"; Period = "Annual" }
)

foreach ($p in $newPeriods) {
    $r = $p.Row
    $compiled.Range("A$r").Value = $p.Label
    $compiled.Range("B$r").Value = $p.Year
    $compiled.Range("C$r").Value = $p.Period
    $compiled.Range("D$r").Value = "2015 CPH"
    $compiled.Range("E$r").Value = "p"
}

# ---------------------------------------------------------------------------
# 2. "Compiled" sheet - corrections to existing rows.
# ---------------------------------------------------------------------------

# Row 247 ("2021 Annual") incorrectly had its Year column set to 2022; fix it to 2021.
$compiled.Range("B247").Value = "2021"

# Row 248 ("2022 Jan") values updated with finalized (non-preliminary) figures, so the
# preliminary "p" flag in column E is cleared.
$compiled.Range("E248").ClearContents()
$compiled.Range("J248").Value = 76347.826000000001
$compiled.Range("K248").Value = 46218.892
$compiled.Range("L248").Value = 43266.404000000002
$compiled.Range("M248").Value = 2952.4879999999998
$compiled.Range("O248").Value = 6430.7969999999996

# Rows 257:259 ("2022 Oct" .. "2022 Dec") - re-share the F/G ratio formulas.
$compiled.Range("F257:G259").Formula = "=K257/J257"

# ---------------------------------------------------------------------------
# 3. "Compiled" sheet - fill in the full figures for the new "2023 Jan" row.
# ---------------------------------------------------------------------------

$compiled.Range("F261").Formula = "=K261/J261"
$compiled.Range("F261").NumberFormat = "0.000"
$compiled.Range("G261").Formula = "=L261/K261"
$compiled.Range("G261").NumberFormat = "0.000"
$compiled.Range("H261").Formula = "=M261/K261"
$compiled.Range("H261").NumberFormat = "0.000"
$compiled.Range("I261").Formula = "=O261/L261"
$compiled.Range("I261").NumberFormat = "0.000"

$compiled.Range("J261").Value = 77104.574999999997
$compiled.Range("J261").NumberFormat = "#,##0"
$compiled.Range("K261").Value = 49724.432999999997
$compiled.Range("K261").NumberFormat = "#,##0"
$compiled.Range("L261").Value = 47351.565000000002
$compiled.Range("L261").NumberFormat = "#,##0"
$compiled.Range("M261").Value = 2372.8690000000001
$compiled.Range("M261").NumberFormat = "#,##0"
$compiled.Range("N261").Formula = "=J261-K261"
$compiled.Range("N261").NumberFormat = "#,##0"
$compiled.Range("O261").Value = 6654.4350000000004
$compiled.Range("O261").NumberFormat = "#,##0"

$compiled.Range("P261").Value = "https://psa.gov.ph/statistics/survey/labor-and-employment/labor-force-survey/title/Unemployment%20Rate%20in%20January%202023%20is%20Estimated%20at%204.8%20Percent"

# ---------------------------------------------------------------------------
# 4. "LFS Reports" sheet - mirrors the new periods added to "Compiled".
# ---------------------------------------------------------------------------

$lfsLabels = @(
    @{ Row = 151; Label = "2022 Annual" },
    @{ Row = 152; Label = "2023 Jan" },
    @{ Row = 153; Label = "2023 Feb" },
    @{ Row = 154; Label = "2023 Mar" },
    @{ Row = 155; Label = "2023 Apr" },
    @{ Row = 156; Label = "2023 May" },
    @{ Row = 157; Label = "2023 Jun" },
    @{ Row = 158; Label = "2023 Jul" },
    @{ Row = 159; Label = "2023 Aug" },
    @{ Row = 160; Label = "2023 Sep" },
    @{ Row = 161; Label = "2023 Oct" },
    @{ Row = 162; Label = "2023 Nov" },
    @{ Row = 163; Label = "2023 Dec" },
    @{ Row = 164; Label = "2023 Annual" }
)

foreach ($l in $lfsLabels) {
    $lfsReports.Range("A$($l.Row)").Value = $l.Label
}

$lfsReports.Range("B152").Value = "https://psa.gov.ph/statistics/survey/labor-and-employment/labor-force-survey/title/Unemployment%20Rate%20in%20January%202023%20is%20Estimated%20at%204.8%20Percent"
